# Update the dSF column (F) values for the wheeler_zack.xlsx sheet.
# New values per row, as derived from the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    3  = -2
    4  = -1
    5  = 3
    7  = 13
    8  = 2
    9  = -1
    10 = -1
    11 = 7
    12 = -2
    13 = 6
    14 = -5
    16 = 8
    17 = 1
    18 = 3
    19 = 9
    20 = -6
    21 = 3
    22 = 2
    24 = 2
    25 = 2
    26 = -3
    27 = 4
    28 = 1
    29 = 5
    30 = 9
    31 = 1
    32 = -3
    33 = -3
    34 = 2
    35 = 1
    36 = 2
    38 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
